$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" -------------------
# This status string shows up on the Overview sheet (columns E "zh-cn" and
# F "de-de") as well as on each language sheet's "Status" column (C).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = "In Translation"

# --- Column width change -----------------------------------------------------
# The status text got shorter ("Ready for handoff" -> "In Translation"), so
# the columns that held it were re-sized (narrower) to fit the new content on
# the Overview sheet (E & F) and on each language sheet (C).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
